$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 3616.7222
$ws.Range("I19").Value = 1849.9166
$ws.Range("J19").Value = 7150.3335
$ws.Range("K19").Value = 1849.9166
$ws.Range("L19").Value = 7150.3335
$ws.Range("M19").Value = -1674.9166
$ws.Range("N19").Value = -7500.3335

$ws.Range("H70").Value = 1499.875
$ws.Range("I70").Value = 999.6667
$ws.Range("J70").Value = 1800
$ws.Range("K70").Value = 2999.0001
$ws.Range("L70").Value = 5400
$ws.Range("M70").Value = -2729.0001
$ws.Range("N70").Value = -5940

$ws.Range("H73").Value = 1499.875
$ws.Range("I73").Value = 999.6667
$ws.Range("J73").Value = 1800
$ws.Range("K73").Value = 2999.0001
$ws.Range("L73").Value = 5400
$ws.Range("M73").Value = -2063.0001
$ws.Range("N73").Value = -7272

$ws.Range("H137").Value = 11139635
$ws.Range("I137").Value = 838.8570999999999
$ws.Range("J137").Value = 50125424
$ws.Range("K137").Value = 2516.5713
$ws.Range("L137").Value = 150376272
$ws.Range("M137").Value = 33.42870000000039
$ws.Range("N137").Value = -150381372

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14118.275
$ws.Range("I32").Value = 12380.607
$ws.Range("J32").Value = 19041.666
$ws.Range("K32").Value = 12380.607
$ws.Range("L32").Value = 19041.666
$ws.Range("M32").Value = -12093.607
$ws.Range("N32").Value = -19615.666

$ws.Range("H117").Value = 27957.084
$ws.Range("J117").Value = 27957.084
$ws.Range("L117").Value = 27957.084
$ws.Range("N117").Value = -37135.084

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H57").Value = 79980
$ws.Range("J57").Value = 79980
$ws.Range("L57").Value = 79980
$ws.Range("N57").Value = -81420

$ws.Range("H82").Value = 18375.354
$ws.Range("I82").Value = 9401.857
$ws.Range("J82").Value = 24656.8
$ws.Range("K82").Value = 9401.857
$ws.Range("L82").Value = 24656.8
$ws.Range("M82").Value = -9018.857
$ws.Range("N82").Value = -25422.8

$ws.Range("H85").Value = 18375.354
$ws.Range("I85").Value = 9401.857
$ws.Range("J85").Value = 24656.8
$ws.Range("K85").Value = 9401.857
$ws.Range("L85").Value = 24656.8
$ws.Range("M85").Value = -8075.857
$ws.Range("N85").Value = -27308.8

$ws.Range("H118").Value = 7755.3193
$ws.Range("J118").Value = 7755.3193
$ws.Range("L118").Value = 7755.3193
$ws.Range("N118").Value = -11069.3193

$ws.Range("H136").Value = 79980
$ws.Range("J136").Value = 79980
$ws.Range("L136").Value = 79980
$ws.Range("N136").Value = -90180

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3747778.5
$ws.Range("I31").Value = 2369.1807
$ws.Range("J31").Value = 55559270
$ws.Range("K31").Value = 2369.1807
$ws.Range("L31").Value = 55559270
$ws.Range("M31").Value = -2074.1807
$ws.Range("N31").Value = -55559860

$ws.Range("H34").Value = 3747778.5
$ws.Range("I34").Value = 2369.1807
$ws.Range("J34").Value = 55559270
$ws.Range("K34").Value = 2369.1807
$ws.Range("L34").Value = 55559270
$ws.Range("M34").Value = -2167.1807
$ws.Range("N34").Value = -55559674

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 46.38889
$ws.Range("I14").Value = 46.38889
$ws.Range("K14").Value = 139.16667
$ws.Range("M14").Value = 33.83332999999999

$ws.Range("H107").Value = 863.85
$ws.Range("I107").Value = 264.2857
$ws.Range("J107").Value = 1186.6923
$ws.Range("K107").Value = 792.8571000000001
$ws.Range("L107").Value = 3560.0769
$ws.Range("M107").Value = 1127.1429
$ws.Range("N107").Value = -7400.0769

$ws.Range("H125").Value = 6329.222
$ws.Range("I125").Value = 4030
$ws.Range("J125").Value = 6616.625
$ws.Range("K125").Value = 12090
$ws.Range("L125").Value = 19849.875
$ws.Range("M125").Value = -7170
$ws.Range("N125").Value = -29689.875

$ws.Range("H132").Value = 2552
$ws.Range("I132").Value = 1062.4
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 9561.6
$ws.Range("L132").Value = 90000
$ws.Range("M132").Value = -7031.6
$ws.Range("N132").Value = -95060

$ws.Range("H134").Value = 4166
$ws.Range("I134").Value = 1909.2307
$ws.Range("J134").Value = 8357.143
$ws.Range("K134").Value = 5727.6921
$ws.Range("L134").Value = 25071.429
$ws.Range("M134").Value = -657.6921000000002
$ws.Range("N134").Value = -35211.429

$ws.Range("H136").Value = 4197.0625
$ws.Range("I136").Value = 924
$ws.Range("J136").Value = 5684.8184
$ws.Range("K136").Value = 2772
$ws.Range("L136").Value = 17054.4552
$ws.Range("M136").Value = 2328
$ws.Range("N136").Value = -27254.4552

$ws.Range("H137").Value = 7179.2856
$ws.Range("I137").Value = 4251.25
$ws.Range("J137").Value = 11083.333
$ws.Range("K137").Value = 12753.75
$ws.Range("L137").Value = 33249.999
$ws.Range("M137").Value = -7653.75
$ws.Range("N137").Value = -43449.999

$ws.Range("H139").Value = 2135.6086
$ws.Range("I139").Value = 1323.591
$ws.Range("J139").Value = 20000
$ws.Range("K139").Value = 3970.773
$ws.Range("L139").Value = 60000
$ws.Range("M139").Value = 1169.227
$ws.Range("N139").Value = -70280

$ws.Range("H140").Value = 3508.9285
$ws.Range("I140").Value = 1284.091
$ws.Range("K140").Value = 3852.273
$ws.Range("M140").Value = 1327.727

$ws.Range("H141").Value = 14300
$ws.Range("I141").Value = 7266.6665
$ws.Range("J141").Value = 21333.334
$ws.Range("K141").Value = 21799.9995
$ws.Range("L141").Value = 64000.00199999999
$ws.Range("M141").Value = -16619.9995
$ws.Range("N141").Value = -74360.00199999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2770.6667
$ws.Range("I126").Value = 1180.5
$ws.Range("J126").Value = 4042.8
$ws.Range("K126").Value = 3541.5
$ws.Range("L126").Value = 12128.4
$ws.Range("M126").Value = -1071.5
$ws.Range("N126").Value = -17068.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 759
$ws.Range("I22").Value = 406.75
$ws.Range("J22").Value = 970.35
$ws.Range("K22").Value = 406.75
$ws.Range("L22").Value = 970.35
$ws.Range("M22").Value = -111.75
$ws.Range("N22").Value = -1560.35

$ws.Range("H27").Value = 759
$ws.Range("I27").Value = 406.75
$ws.Range("J27").Value = 970.35
$ws.Range("K27").Value = 406.75
$ws.Range("L27").Value = 970.35
$ws.Range("M27").Value = -299.75
$ws.Range("N27").Value = -1184.35

$ws.Range("H132").Value = 12827127
$ws.Range("I132").Value = 3217.0908
$ws.Range("K132").Value = 9651.2724
$ws.Range("M132").Value = -7121.2724
